$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row D/E value updates ---
$ws.Range("D2").Value = "61.487.35"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "3.445.93"
$ws.Range("E3").Value = "  +2.13%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "578.49"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").Value = "143.95"
$ws.Range("E6").Value = "  +5.43%  "

$ws.Range("D7").Value = "3.447.64"
$ws.Range("E7").Value = "  +2.27%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("D10").Value = "7.61"
$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("E11").Value = "  +3.00%  "

$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").Value = "4.034.84"
$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").Value = "27.93"
$ws.Range("E14").Value = "  +9.05%  "

$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D18").Value = "61.639.36"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  +9.31%  "

$ws.Range("D20").Value = "14.26"
$ws.Range("E20").Value = "  +3.59%  "

$ws.Range("D21").Value = "9.51"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("D22").Value = "388.59"
$ws.Range("E22").Value = "  +3.82%  "

$ws.Range("D23").Value = "0.563"
$ws.Range("E23").Value = "  +3.17%  "

$ws.Range("D24").Value = "73.37"
$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("D25").Value = "5.77"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").Value = "3.590.09"
$ws.Range("E28").Value = "  +2.11%  "

$ws.Range("E29").Value = "  +1.46%  "

$ws.Range("D30").Value = "7.59"
$ws.Range("E30").Value = "  +3.43%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +1.53%  "

$ws.Range("D33").Value = "2.18"

$ws.Range("E34").Value = "  -11.32%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").Value = "23.98"
$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("D37").Value = "3.476.31"
$ws.Range("E37").Value = "  +2.42%  "

$ws.Range("D38").Value = "6.99"
$ws.Range("E38").Value = "  +3.26%  "

$ws.Range("D39").Value = "5.11"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").Value = "1.55"
$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").Value = "166.66"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("D42").Value = "28.08"
$ws.Range("E42").Value = "  +13.60%  "

$ws.Range("D43").Value = "0.0781"
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("D44").Value = "0.803"
$ws.Range("E44").Value = "  +3.73%  "

$ws.Range("D47").Value = "4.47"
$ws.Range("E47").Value = "  +3.76%  "

$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").Value = "2.579.35"
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  -1.65%  "

$ws.Range("D51").Value = "6.92"
$ws.Range("E51").Value = "  +2.32%  "
# --- Row 16/17 swap: ShibaInu/WrappedEther order flip with new values ---
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.463.00"
$ws.Range("E16").Value = "  +2.58%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +1.57%  "

# --- Row 45/46 swap: OKB/FirstDigitalUSD order flip with new values ---
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "42.33"
$ws.Range("E46").Value = "  +1.52%  "
